# Re-run of the Data Preprocessing & Cleaning script ("enriched_data.xlsx").
# This replays a fresh pass of the enrichment pipeline: a handful of rows get
# corrected "subregion"/"capital" lookups, a handful of rows get refreshed
# population/area figures (which ripple into the derived population_density
# and language_density columns), and every row's "timestamp" column is
# stamped with the time of this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Corrected lookup values (subregion / capital) picked up by the refresh.
# ---------------------------------------------------------------------------

# Sierra Leone (row 5): capital lookup resolved.
$ws.Range("J5").Value = '["Freetown"]'

# Barbados (row 9): subregion lookup now reports Unknown.
$ws.Range("F9").Value = "Unknown"

# Benin (row 14): capital lookup resolved.
$ws.Range("J14").Value = '["Porto-Novo"]'

# France (row 24): subregion lookup resolved.
$ws.Range("F24").Value = "Western Europe"

# Mayotte (row 81): capital lookup now reports Unknown.
$ws.Range("J81").Value = "Unknown"

# Malta (row 214): capital lookup now reports Unknown.
$ws.Range("J214").Value = "Unknown"

# Kuwait (row 219): subregion lookup resolved.
$ws.Range("F219").Value = "Western Asia"

# Saint Martin (row 242): subregion lookup now reports Unknown.
$ws.Range("F242").Value = "Unknown"

# ---------------------------------------------------------------------------
# 2) Refreshed population / area figures, with their dependent derived
#    columns (population_density = population / area;
#    language_density = language_count * 1e6 / population) recomputed.
# ---------------------------------------------------------------------------

# Aruba (row 47): population refreshed.
$ws.Range("G47").Value = 109436
$ws.Range("O47").Value = 607.9777777777778
$ws.Range("U47").Value = 18.27552176614642

# Antarctica (row 68): population refreshed.
$ws.Range("G68").Value = 1000
$ws.Range("O68").Value = 0.00007142857142857143

# Pakistan (row 76): area refreshed.
$ws.Range("H76").Value = 881912
$ws.Range("O76").Value = 250.4698099130072

# Bosnia and Herzegovina (row 113): area refreshed.
$ws.Range("H113").Value = 51209
$ws.Range("O113").Value = 64.06715616395556

# Western Sahara (row 130): population refreshed.
$ws.Range("G130").Value = 510713
$ws.Range("O130").Value = 1.919973684210526
$ws.Range("U130").Value = 5.874140662172296

# Thailand (row 134): area refreshed.
$ws.Range("H134").Value = 513102.2526792237
$ws.Range("O134").Value = 136.0352203396715

# Tajikistan (row 213): population refreshed.
$ws.Range("G213").Value = 9294356
$ws.Range("O213").Value = 64.95007686932215
$ws.Range("U213").Value = 0.2151843548923669

# Zimbabwe (row 241): area refreshed.
$ws.Range("H241").Value = 388670.346818461
$ws.Range("O241").Value = 38.24044494688999

# ---------------------------------------------------------------------------
# 3) Stamp every data row's "timestamp" column (N) with this run's time.
#    The pipeline writes rows in small batches, so the timestamp advances a
#    few times across the sheet rather than being a single constant value.
# ---------------------------------------------------------------------------

$ws.Range("N2:N10").Value = "2025-03-31T00:48:34.412194"
$ws.Range("N11:N20").Value = "2025-03-31T00:48:34.427804"
$ws.Range("N21:N29").Value = "2025-03-31T00:48:34.443426"
$ws.Range("N30:N38").Value = "2025-03-31T00:48:34.459052"
$ws.Range("N39:N47").Value = "2025-03-31T00:48:34.474676"
$ws.Range("N48:N56").Value = "2025-03-31T00:48:34.490306"
$ws.Range("N57:N66").Value = "2025-03-31T00:48:34.505932"
$ws.Range("N67:N75").Value = "2025-03-31T00:48:34.522110"
$ws.Range("N76:N84").Value = "2025-03-31T00:48:34.537143"
$ws.Range("N85:N90").Value = "2025-03-31T00:48:34.552812"
$ws.Range("N91:N99").Value = "2025-03-31T00:48:34.568435"
$ws.Range("N100:N108").Value = "2025-03-31T00:48:34.584062"
$ws.Range("N109:N117").Value = "2025-03-31T00:48:34.599688"
$ws.Range("N118:N126").Value = "2025-03-31T00:48:34.615316"
$ws.Range("N127:N135").Value = "2025-03-31T00:48:34.630937"
$ws.Range("N136:N145").Value = "2025-03-31T00:48:34.646569"
$ws.Range("N146:N154").Value = "2025-03-31T00:48:34.662194"
$ws.Range("N155:N163").Value = "2025-03-31T00:48:34.677812"
$ws.Range("N164:N172").Value = "2025-03-31T00:48:34.693440"
$ws.Range("N173:N182").Value = "2025-03-31T00:48:34.709060"
$ws.Range("N183:N191").Value = "2025-03-31T00:48:34.724697"
$ws.Range("N192:N200").Value = "2025-03-31T00:48:34.740314"
$ws.Range("N201:N209").Value = "2025-03-31T00:48:34.755943"
$ws.Range("N210:N218").Value = "2025-03-31T00:48:34.771570"
$ws.Range("N219:N227").Value = "2025-03-31T00:48:34.787194"
$ws.Range("N228:N236").Value = "2025-03-31T00:48:34.802888"
$ws.Range("N237:N245").Value = "2025-03-31T00:48:34.818443"
$ws.Range("N246:N251").Value = "2025-03-31T00:48:34.834067"
